$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 80318
$ws.Range("B2").Value = "Heitor Santos"
$ws.Range("C2").Value = "Operacoes"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45085
$ws.Range("G2").Value = 4985.5

# Row 3
$ws.Range("A3").Value = 14563
$ws.Range("B3").Value = "José Pedro Fernandes"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Doenca"
$ws.Range("F3").Value = 45097
$ws.Range("G3").Value = 4801.01

# Row 4
$ws.Range("A4").Value = 27436
$ws.Range("B4").Value = "Bianca da Cunha"
$ws.Range("C4").Value = "Atendimento ao Cliente"
$ws.Range("D4").Value = "Viagem de negocios"
$ws.Range("G4").Value = 5220.23

# Row 5
$ws.Range("A5").Value = 78828
$ws.Range("B5").Value = "Pietra da Costa"
$ws.Range("C5").Value = "Operacoes"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45090
$ws.Range("G5").Value = 8116.24

# Row 6
$ws.Range("A6").Value = 46482
$ws.Range("B6").Value = "Fernando Santos"
$ws.Range("C6").Value = "Marketing"
$ws.Range("D6").Value = "Doenca"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45087
$ws.Range("G6").Value = 5744.2

# Row 7
$ws.Range("A7").Value = 90911
$ws.Range("B7").Value = "Catarina Castro"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Doenca"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45093
$ws.Range("G7").Value = 2203.87

# Row 8
$ws.Range("A8").Value = 13253
$ws.Range("B8").Value = "Alícia Nogueira"
$ws.Range("C8").Value = "Juridico"
$ws.Range("D8").Value = "Viagem de negocios"
$ws.Range("F8").Value = 45085
$ws.Range("G8").Value = 9902.43

# Row 9
$ws.Range("A9").Value = 31123
$ws.Range("B9").Value = "Srta. Maria Julia Azevedo"
$ws.Range("C9").Value = "Atendimento ao Cliente"
$ws.Range("D9").Value = "Consulta medica"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45094
$ws.Range("G9").Value = 9752.450000000001

# Row 10
$ws.Range("A10").Value = 49999
$ws.Range("B10").Value = "Marina Cirino"
$ws.Range("C10").Value = "Juridico"
$ws.Range("D10").Value = "Doenca"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45091
$ws.Range("G10").Value = 3241.59

# Row 11
$ws.Range("A11").Value = 11969
$ws.Range("B11").Value = "Eloá Moreira"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45088
$ws.Range("G11").Value = 8281.219999999999
